$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emDash = [char]0x2014

$ws.Range("C36").Value = "[name=`"Big Bob`"]   '...Monster?'`n"
$ws.Range("C68").Value = "[name='Captain']   Now! Fiiire!`n"
$ws.Range("C70").Value = "[name='Captain']   That's her right there, the walking Catastrophe! Hurry, shoot her down!`n"
$ws.Range("C72").Value = "[name='Captain']   What do you think you're doing, kid?`n"
$ws.Range("C74").Value = "[name='Captain']   G-get out of my way, you rascal$emDash!`n"
